$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Set B5 to a single space value (shown in diff as new shared string " ")
$ws.Range("B5").Value = " "

# Reflect the selection change that Excel records when the user selects B5
$ws.Range("B5").Select()
